$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = '[''이벤트/쿠폰 > Little SSG 어린이날 선물대전'', ''스마일클럽'', ''20% 쿠폰'', ''멤버십 전용 유아동 20% 쿠폰'', ''선착순 4만장 SmileClub전용 상품쿠폰 20%'', ''3만원 이상 상품 구매시 최대 7천원 할인, 일부 특가상품 적용 제외, 신세계백화점, 신세계몰, 이마트몰'', ''쿠폰은 25일 10시에 오픈됩니다.'', ''선착순 쿠폰 마감 되었습니다. 감사합니다'', ''발급 후 [My SSG] > [쿠폰] > [보유쿠폰] 에서 확인 가능합니다.'', ''[ 미미월드 ] 미리준비하는 어린이날 미미/ 똘똘이/ 뽀로로 외 최대 30%할인'', ''어린이날 실바니안패밀리 브레드이발소 산리오 외 세트특가할인'', ''★어린이날★ 산리오/몰랑 완구 모음전 최대 30%할인'', ''최대 30% 특가'', ''[해즈브로] 플레이도 3만원 이상 구매시 매트 증정 *소진시 행사종료*'', ''할인이 쓱쓱 인기 여아캐릭터 완구 모음전 미미월드 토이트론 외'', ''★쓱세일★ 해즈브로 인기 브랜드 할인모음'', ''시크릿쥬쥬/ LOL/ 미미 최대 30% 할인'', ''[쁘띠엘린] 완구 특가 립프로그/젤리캣/하퍼스테이블 할인'', ''[꼬미마녀라라 외] 인기 신상완구 특가 할인!'']'
$ws.Range("A7").Value = '바이오퍼블릭 기프트 위크'
$ws.Range("B7").Value = '15% 쿠폰으로 건강을 선물하세요!'
$ws.Range("C7").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000005633'
$ws.Range("D7").Value = '바이오퍼블릭 기프트 위크'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2023.04.24'
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = '2023.04.30'
$ws.Range("G7").Value = '[''이벤트/쿠폰 > 바이오퍼블릭 기프트 위크'', ''스마일클럽'', ''+카드청구할인'', ''~15%쿠폰'', ''50%~, 1+1'', ''건강식품 ~15%쿠폰'', ''카드청구할인까지'', ''선착순  1천명'', ''장바구니쿠폰'', ''%'', ''5만원 이상 구매 시 최대 1만 5천원 할인'', ''쿠폰 사용 기간: 발급 당일'', ''3만원 이상 구매 시 최대 1만 5천원 할인'', "쿠폰 발급 후''My SSG ▶ 쿠폰 ▶ 보유쿠폰'' 에서 확인 가능합니다.", ''[SSGPAY전용] 신한카드 8만원 이상 5% 청구할인 (일 10만원 한) 자세히보기'', ''[SSGPAY전용] 삼성카드 8만원 이상 5% 청구할인 (일 10만원 한) 자세히보기'', ''SSG.COM 삼성카드 [SSGPAY전용] 8만원 이상 7% 청구할인 (일 10만원 한) 자세히보기'', ''국민카드 8만원 이상 5% 청구할인 (일 10만원 한) 자세히보기'', ''비씨카드[SSGPAY전용] 8만원 이상 5% 청구할인 (일 10만원 한) 자세히보기'', ''~15%쿠폰으로 건강 선물하기'', ''[오쏘몰] 프리미엄 종합 비타민 특별할인 무료배송 및 ~17% OFF!'', ''[Biopublic]100% 식물성 프로틴 케어 350g(25g*14포)'', ''[이마트해외직구] 신타6 아이솔레이트 단백질파우더 1.82kg 1+1 모음 , 무료배송 1+1'', ''50%~,  1+1 득템 찬스'', ''50%~ 할인 상품'', ''50% 할인'', ''[함소아제약 외]건강식품 50%이상 할인된 가격'', ''[기프트위크](50%할인) 휴럼 건강식품 반값행사'', ''홍삼정스틱네이처 10mlx28 1세트+쇼핑백 外 50% 할인특가'', ''[기프트위크][쓱배송]건강식품 반값행사'', ''50% 할인. 콤비타'', ''[쏜리서치] 공식판매처! 인기 영양제 1+1 구성 특가! 2/데이,S.A.T 외 무료배송'', ''[2개 세트][무료배송] 오쏘몰 Orthomol 이뮨 드링크 + 정제 30일분'', ''[독일직구][무료배송] 오쏘몰 이뮨 30일 1+1 2박스 드링크 정제 멀티비타민 orthomol / 최신제조일자 / 빠른배송 / 독일현지배송'', ''[최대40%할인] 건강식품 선물하기'', ''[쓱배송]락티브 비타민 가격할인'', ''어린이 멀티비타민 유산균 오메가3 백화점 선물세트 [사은품증정이벤트]'', ''트리플러스50+/90캡슐 x 3개 [사은품증정이벤트]'', ''[쓱배송]루테인 초특가'', ''[쓱배송][정관장]홍삼진본 40ml*20포[쇼핑백포함]'', ''[15% 즉시할인][세노비스] 키즈 멀티비타민미네랄 (1.5g*90정*2통/총90일분) ★미니 멀티비타민 2개 증정★'', ''[15% 즉시할인][세노비스] 트리플러스 맨 (90캡슐) + 트리플러스 우먼 (90캡슐) (커플세트) ★미니 트리플러스 맨+우먼 증정★'', ''★139,000→69,500[쓱배송] 휴럼 천년침향환(3.7g*60환)[쇼핑백 동봉]'', ''★ 레이델 특가상품 모음★ 레이델 폴리코사놀/비즈왁스알코올 셀100  外~ 사은품 + 쿠폰☆'', ''[쓱배송][종근당건강] 6년근홍삼정업 2개입세트(240g*2병)[쇼핑백 동봉]'', ''[무료배송][오쏘몰] 이뮨(immun) 드링크+정제 30일분 3개묶음 / 독일직배송 / 최신제조일자 / 본사 공식'', ''최대 20%쿠폰  + 청구할인'', ''~15% 쿠폰가능 + 청구할인'', ''10%할인 + 쇼핑백증정'', ''최대 54% 할인'', ''30% + 5% SSG머니 할인'', ''1+1특가 + 추가 4% 즉시할인'', ''2+1 증정 이벤트'']'
$ws.Range("A8").Value = 'SSG.COM 앰배서더 SSG 메이트 모집'
$ws.Range("B8").Value = 'SSG 메이트가 되어 믿고 사는 즐거움을 누려보세요'
$ws.Range("C8").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000005708'
$ws.Range("D8").Value = '[SSG.COM 앰배서더] SSG 메이트 모집'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2023.04.27'
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = '2023.05.10'
$ws.Range("G8").Value = '[''이벤트/쿠폰 > [SSG.COM 앰배서더] SSG 메이트 모집'', ''스마일클럽'']'
$ws.Range("G10").Value = '[''이벤트/쿠폰 > SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'', ''스마일클럽'', ''SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'', ''SSG머니 최대 15% 적립'', ''스마일클럽 월 이용료 3,900원 결제일 할인'', ''최대 15% SSG머니 적립 중 5%는 SSG.COM에서 제공하는 서비스로 자세한 내용은 SSG.COM 확인바람.'', ''01. SSG MONEY 최대 15% 적립!'', ''       이용실적 관계없이 적립 : 1 % + 이용실적 충족 시 적립 4% = 적립 가능한 최대 SSG머니 5%'', ''       카드혜택 + 스마일클럽 혜택'', ''       이용실적 관계없이 적립 1% + 이용실적 충족 시 적립 9% + 쓱배송/쓱배송 traders/새벽배송 상품 구매 시 5% = 적립 가능한 최대 SSG머니 15%'', ''최대 15% SSG MONEY 적립'', ''SSG.COM 삼성카드 최대 10% 적립+스마일클럽 5% 적립'', ''SSG.COM 삼성카드 최대 10% 적립(1% 적립+9% 추가 적립)'', ''1% 적립(전월 실적 조건 및 적립한도 없음)'', ''9% 추가 적립'', ''삼성카드 할인이 적용된 일시불 및 할부 이용금액은 제외됩니다.'', ''스마일클럽 5% 적립'', ''쓱배송/쓱배송 traders/새벽배송 상품 결제건에 한해 혜택을 받을 수 있습니다.'', ''02. 스마일클럽 월이용료 매월 3,900원 결제일 할인'', ''15,000원만 담아도 쓱 새벽배송 무료배송(SSG가입시)'', ''1쓱 새벽배송 트레이더 - SSG머니 최대 5% 적립'', ''장바구니 최대 10% 할인쿠폰'', ''스마일클럽 단독 혜택 - 스타벅스 상품 전용 딜'', ''매월 4장씩 최대 12% 할인쿠폰'', ''스마일배송 1만 5천원 이상 무제한 무료배송'', ''스마일배송 상품 스마일캐시 1% 적립'', ''카드 할인 혜택 자세히보기 (레이어팝업 열기)'', ''2022.10.26부터, SSG.COM 삼성카드 발급 시 스마일클럽에 자동 가입됩니다.'', ''G마켓 또는 옥션을 통해 스마일클럽에 가입한 경우 혜택을 받을 수 없습니다.'', ''SSG.COM 삼성카드로 스마일클럽 월 이용료(3,900원) 정기결제 시 혜택이 제공됩니다.(월 1회)'', ''SSG.COM을 통한 스마일클럽 가입건에 한해 혜택이 제공됩니다.'', ''결제금액이 할인금액보다 적을 경우, 결제금액만큼 할인이 적용됩니다.(결제금액이 없는 경우 할인 대상에서 제외)'', ''5만 6천원 이상 결제 시 사용 가능한 5만 5천원 할인쿠폰 제공'', ''55,000원 할인쿠폰 - SSGPAY 바로결제 이용 시 사용 가능'', ''쿠폰 발급기간 : 2023.04.01 ~ 2023.04.30'', ''쿠폰 사용기간 : 2023.04.01 ~ 2023.04.30'', ''SSGPAY 바로결제 이용 이력이 없는 회원'', ''SSGPAY 바로결제에 등록된 SSG.COM 삼성카드로 결제해야 쿠폰을 사용할 수 있습니다.'', ''배송비 등을 제외한 최종 결제금액이 56,000원 이상이여야 쿠폰을 사용할 수 있습니다.'', ''SSGPAY 바로결제에 등록된 SSG.COM 삼성카드로 결제하셔야 쿠폰을 사용할 수 있습니다.'', ''쿠폰은 통합 회원 본인 명의의 SSG.COM 삼성카드로 결제 시 사용 가능하며, 다른 부정적인 방법으로 사용한 경우에는 주문이 취소될 수 있습니다.'', ''쿠폰은 결제 화면에서 자동으로 적용됩니다.'', ''결제 화면에서 쿠폰 변경을 원할 경우 ‘쿠폰선택’을 눌러주세요.'', ''본 상품 발급 및 SSGPAY앱 설치시 바로결제에 자동 등록 됩니다.'', ''SSG MONEY는 매월 1일~말일까지 매출전표가 접수된 금액에 대해 다음달 25일 SSG.COM 계정으로 적립됩니다.'', ''적립된 SSG MONEY는 SSGPAY 회원가입 후 조회 및 사용 가능합니다. 단, SSG.COM 에서는 SSG.COM만 가입해도 사용 가능합니다.'', ''연체이자율 : 회원별/이용상품별 정상이자율+3.0%p(최고 연 20.0%)'', ''이미 SSG.COM 삼성카드를 가지고 계시네요!스마일클럽 가입하고 모든 혜택 누리세요'', ''이미 SSG.COM 삼성카드를 가지고 계시네요!G마켓 또는 옥션을 통해 스마일클럽에 가입한 경우, 월 이용료 결제일할인 혜택을 받을 수 없습니다.'', ''쿠폰 다운되었습니다.다운받은 쿠폰은 SSG.COM 삼성카드 발급 후 바로 사용 가능합니다.'', ''이벤트 쿠폰이 이미 발급되었습니다.다운받은 쿠폰은 SSG.COM 삼성카드 발급 후 바로 사용 가능합니다.'']'
$ws.Range("A11").Value = '다이슨 에어랩 컴플리트 4/27(목) 7PM'
$ws.Range("B11").Value = '롱배럴(20/40mm) 추가 증정+7만원 적립! 인기컬러 구매 전.고.객. 모두 증정!'
$ws.Range("C11").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000005778'
$ws.Range("D11").Value = '다이슨@SSG.LIVE 4/27(수) 7:00PM'
$ws.Range("G11").Value = '[''이벤트/쿠폰 > 다이슨@SSG.LIVE 4/27(수) 7:00PM'', ''스마일클럽'', ''- 사은품 지급 및 이벤트 혜택 제공'', ''- 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기'']'
$ws.Range("A12").Value = '사봉 4/27(목) 8PM'
$ws.Range("B12").Value = '방송중에만! 원더오브자스민外 연중 최대할인 & 전구매고객 단독사은품 증정'
$ws.Range("C12").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000005664&domainSiteNo=6005'
$ws.Range("D12").Value = '사봉 @SSG.LIVE 4/27(목) 20:00PM'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2023.04.17'
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = '2023.04.27'
$ws.Range("G12").Value = '[''이벤트/쿠폰 > 사봉 @SSG.LIVE 4/27(목) 20:00PM'', ''스마일클럽'', ''                     SSG.LIVE 사은품 지급 및 이벤트 혜택 당첨 주의사항'', '' - 사은품 지급 및 이벤트 혜택 제공'', '' - 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기'']'
$ws.Range("A13").Value = '[공동판매구역] 맛젤 4/28(금) 11AM'
$ws.Range("B13").Value = '제스프리 썬골드키위 시즌 컴백! 핫딜가 19,900원 & 구매 시 스푼 100% 증정'
$ws.Range("C13").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000005691&domainSiteNo=6005'
$ws.Range("D13").Value = '[SSG.LIVE]4/28(금) 11AM 제스프리 골드키위'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2023.04.20'
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = '2023.04.28'
$ws.Range("G13").Value = '[''이벤트/쿠폰 > [SSG.LIVE]4/28(금) 11AM 제스프리 골드키위'', ''스마일클럽'', ''- 사은품 지급 및 이벤트 혜택 제공'', ''- 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기'']'
$ws.Range("G14").Value = '[''이벤트/쿠폰 > [4/24~30 ] 럭셔리 더 뷰티 가정의 달 GIFT'', ''스마일클럽'', ''15%'', ''10%'', ''기획상품10% 할인 전구매 고객 두피세럼 증정'', ''[시코르 추천] 인기 뷰티템 UP TO 30%'', ''자연의 작품, 연작   카밍 앤 컴포팅라인 20% 할인 + 쓱단독 사은품 증정'', ''[몰튼 브라운]SSG 단독 선물포장+UP TO 35%'']'
